# first cut at adding ibgp mesh betwwen pes
$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet ("core_interfaces+")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "vpn_ibgp+"

# Header row
$ws.Range("A1").Value = "host"
$ws.Range("B1").Value = "@group"
$ws.Range("C1").Value = "neighbors+.name"
$ws.Range("A1:C1").Style = "header"

# iBGP full mesh between PE loopbacks
$pes = @(
    @("host_vars/nyc-1.yaml", "10.52.100.3"),
    @("host_vars/nyc-2.yaml", "10.52.100.4"),
    @("host_vars/bos-1.yaml", "10.52.100.5"),
    @("host_vars/bos-2.yaml", "10.52.100.6"),
    @("host_vars/lax-1.yaml", "10.52.100.11"),
    @("host_vars/lax-2.yaml", "10.52.100.12")
)

$r = 2
foreach ($pe in $pes) {
    foreach ($other in $pes) {
        if ($other[0] -ne $pe[0]) {
            $ws.Cells.Item($r, 1).Value = $pe[0]
            $ws.Cells.Item($r, 2).Value = "VPN_iBGP"
            $ws.Cells.Item($r, 3).Value = $other[1]
            $r = $r + 1
        }
    }
}

$lastRow = $r - 1
$ws.Range("A2:C$lastRow").Style = "value"

# Column widths
$ws.Columns.Item(1).ColumnWidth = 27.5703125
$ws.Columns.Item(2).ColumnWidth = 11.5703125
$ws.Columns.Item(3).ColumnWidth = 17

# Footer to match the rest of the workbook
$ws.PageSetup.OddFooter = "&C&1#&""Calibri""&7&K000000Juniper Business Use Only"

# Selection / view state for the new active sheet
$ws.Range("C$lastRow").Select()
